$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 18-19; everything from the old row 18 down
# (old rows 18..49) shifts down to rows 20..51.
$ws.Range("A18:A19").EntireRow.Insert()

# New row 18: Vega Monumental Concepción - Pepino dulce, Primera, 27-May-2022
$d = Get-Date -Year 2022 -Month 5 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(18,1).Value  = 11
$ws.Cells.Item(18,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(18,3).Value  = "Bíobío"
$ws.Cells.Item(18,4).Value  = $d
$ws.Cells.Item(18,5).Value  = 8
$ws.Cells.Item(18,6).Value  = 100112043
$ws.Cells.Item(18,7).Value  = "Pepino dulce"
$ws.Cells.Item(18,8).Value  = "Cultivar IV Región"
$ws.Cells.Item(18,9).Value  = "Primera"
$ws.Cells.Item(18,10).Value = 100
$ws.Cells.Item(18,11).Value = 13000
$ws.Cells.Item(18,12).Value = 14000
$ws.Cells.Item(18,13).Value = 13500
$ws.Cells.Item(18,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(18,15).Value = "Provincia de Limarí"
$ws.Cells.Item(18,16).Value = 750
$ws.Cells.Item(18,17).Value = 18
$ws.Cells.Item(18,18).Value = "Hortaliza"

# New row 19: Vega Monumental Concepción - Pepino dulce, Segunda, 27-May-2022
$ws.Cells.Item(19,1).Value  = 11
$ws.Cells.Item(19,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(19,3).Value  = "Bíobío"
$ws.Cells.Item(19,4).Value  = $d
$ws.Cells.Item(19,5).Value  = 8
$ws.Cells.Item(19,6).Value  = 100112043
$ws.Cells.Item(19,7).Value  = "Pepino dulce"
$ws.Cells.Item(19,8).Value  = "Cultivar IV Región"
$ws.Cells.Item(19,9).Value  = "Segunda"
$ws.Cells.Item(19,10).Value = 50
$ws.Cells.Item(19,11).Value = 11000
$ws.Cells.Item(19,12).Value = 11000
$ws.Cells.Item(19,13).Value = 11000
$ws.Cells.Item(19,14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(19,15).Value = "Provincia de Limarí"
$ws.Cells.Item(19,16).Value = 611
$ws.Cells.Item(19,17).Value = 18
$ws.Cells.Item(19,18).Value = "Hortaliza"
